$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: logistic_embeddings
$ws.Range("C5").Value = 0.18
$ws.Range("D5").Value = 0.327
$ws.Range("E5").Value = 0.383
$ws.Range("F5").Value = 0.411
$ws.Range("G5").Value = 0.46
$ws.Range("H5").Value = 0.475

# Row 7: classical-best-embeddings -> classical-best-embed
$ws.Range("A7").Value = "classical-best-embed"
$ws.Range("C7").Value = 0.18
$ws.Range("D7").Value = 0.327
$ws.Range("E7").Value = 0.383

# Row 8: BERT-base
$ws.Range("C8").Value = 0.169
$ws.Range("D8").Value = 0.426
$ws.Range("E8").Value = 0.47
$ws.Range("F8").Value = 0.546
$ws.Range("G8").Value = 0.591
$ws.Range("H8").Value = 0.596

# Row 9: BERT-base-nli
$ws.Range("B9").Value = 0.382
$ws.Range("C9").Value = 0.472
$ws.Range("D9").Value = 0.532
$ws.Range("E9").Value = 0.56
$ws.Range("F9").Value = 0.611
$ws.Range("G9").Value = 0.63
$ws.Range("H9").Value = 0.639
